{"js": "// The document repeats the campaign-dates blurb in several places; each\n// copy contains a duplicated word: \"...zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed\n// Pegasus...\". The fix removes the stray lower-case \"souhv\u011bzd\u00ed \" so the\n// sentence reads \"...zobrazuj\u00edSouhv\u011bzd\u00ed Pegasus...\".\nconst body = context.document.body;\nconst results = body.search(\"zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"zobrazuj\u00edSouhv\u011bzd\u00ed\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace every occurrence of \"zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed\" with\n# \"zobrazuj\u00edSouhv\u011bzd\u00ed\" throughout the document body (removes the stray\n# duplicated lower-case \"souhv\u011bzd\u00ed\" word and the space before it).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"zobrazuj\u00ed souhv\u011bzd\u00ed Souhv\u011bzd\u00ed\"\n$find.Replacement.Text = \"zobrazuj\u00edSouhv\u011bzd\u00ed\"\n$find.Forward = $true\n$find.Wrap = 2            # wdFindContinue - keep going across the whole story\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute([ref]$find.Text, $false, $true, $false, $false, $false, $true, 2, $false, [ref]$find.Replacement.Text, 2)\n"}
